$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change B3 value from 12 to 16
$ws.Range("B3").Value = 16

# Remove the stray J9 cell (shared string "S") entirely
$ws.Range("J9").ClearContents()

# Add new data rows 16-19
$ws.Range("A16").Value = 6
$ws.Range("B16").Value = 3
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = "default"

$ws.Range("A17").Value = 6
$ws.Range("B17").Value = 4
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = "default"

$ws.Range("A18").Value = 7
$ws.Range("B18").Value = 5
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = "default"

$ws.Range("A19").Value = 7
$ws.Range("B19").Value = 6
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = "default"

# Update the active selection to match the new last-edited cell
$ws.Range("B19").Select()
